# Revert "Added converter option"
# This reverts the earlier commit that renamed the "type" converter option
# (and its "string" value) to "readAs" / "text". Restore the original
# option names/values on the affected cells.

$wb = $excel.ActiveWorkbook

# --- Sheet "limit" ---
$wsLimit = $wb.Worksheets.Item("limit")
$wsLimit.Range("C1").Value = "options?limit=5#aaa?type=string"

# --- Sheet "offset" ---
$wsOffset = $wb.Worksheets.Item("offset")
$wsOffset.Range("C1").Value = "optionsOffset?limit=5&offset=2#aaa?type=string"

# --- Sheet "inTableOptions" ---
$wsInTable = $wb.Worksheets.Item("inTableOptions")
$wsInTable.Range("A3").Value = "single?type"
$wsInTable.Range("C3").Value = "string"
$wsInTable.Range("A8").Value = "optionsInTable?type"
$wsInTable.Range("E8").Value = "string"

# --- Restore the previously-active sheet/selection ---
# Before the reverted commit, "optionForTableAndColumn" was the active tab
# (with C1 selected) and "inTableOptions" was not the selected tab
# (its selection sits on C1 instead of A9).
$wsInTable.Range("C1").Select() | Out-Null

$wsOptionForTableAndColumn = $wb.Worksheets.Item("optionForTableAndColumn")
$wsOptionForTableAndColumn.Activate()
$wsOptionForTableAndColumn.Range("C1").Select() | Out-Null
